# Add "Done Y/N?" = Yes and "Team Member" = Oisin to every task row
# in the WDD project task list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$taskRows = @(4, 5, 6, 7, 8, 11, 12, 13, 14, 15, 16, 19, 24, 26, 29, 31, 32, 33, 34, 37)

foreach ($r in $taskRows) {
    $ws.Cells.Item($r, 3).Value = "Yes"
    $ws.Cells.Item($r, 4).Value = "Oisin"
}

# Move the active selection to the last filled cell, as in the saved file.
$ws.Range("D37").Select()
